$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.037.67"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.849.42"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "234.78"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "42.53"
$ws.Range("E8").Value = "  +7.20%  "
$ws.Range("D9").Value = "0.330"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "0.0984"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").Value = "2.116.29"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.849.11"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "11.35"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "0.676"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "4.69"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "35.008.19"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "70.02"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "240.87"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "12.13"
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").Value = "4.80"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "2.28"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "170.89"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  +20.57%  "
$ws.Range("D27").Value = "7.87"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").Value = "17.64"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "0.0556"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").Value = "3.98"
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("D33").Value = "3.97"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +23.73%  "
$ws.Range("D35").Value = "2.00"
$ws.Range("E35").Value = "  +11.48%  "
$ws.Range("D36").Value = "0.771"
$ws.Range("E36").Value = "  +9.84%  "
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("E38").Value = "  +10.91%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "0.0201"
$ws.Range("E40").Value = "  +4.74%  "
$ws.Range("D41").Value = "1.345.67"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").Value = "14.92"
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").Value = "  +4.98%  "
$ws.Range("D44").Value = "12.86"
$ws.Range("E44").Value = "  +88.49%  "
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").Value = "2.75"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("D49").Value = "2.028.21"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  +16.11%  "
$ws.Range("E51").Value = "  +1.49%  "
